$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A: move ids zero-padded
$ws.Range("A2").Value = "z0bug.move_01"
$ws.Range("A3").Value = "z0bug.move_02"
$ws.Range("A4").Value = "z0bug.move_03"
$ws.Range("A5").Value = "z0bug.move_04"
$ws.Range("A6").Value = "z0bug.move_05"
$ws.Range("A7").Value = "z0bug.move_06"
$ws.Range("A8").Value = "z0bug.move_07"
$ws.Range("A9").Value = "z0bug.move_08"
$ws.Range("A10").Value = "z0bug.move_09"
$ws.Range("A11").Value = "z0bug.move_10"

# Column C: date refs zero-padded
$ws.Range("C2").Value = "<002-12-99"
$ws.Range("C3").Value = "<001-01-15"
$ws.Range("C4").Value = "<001-01-31"
$ws.Range("C5").Value = "<001-04-05"
$ws.Range("C6").Value = "<001-07-05"
$ws.Range("C7").Value = "<001-10-05"

# Column E: journal_id changed for all rows
$ws.Range("E2:E11").Value = "external.MISC"

# Column F: fiscalyear_id changed for all rows
$ws.Range("F2:F11").Value = "z0bug.fy_%(year)s"

# Selection moved to E3:E11
$ws.Range("E3:E11").Select()
